$wb = $excel.ActiveWorkbook

# Disaggregation of commodity "Copper ores and concentrates" -> "Copper".
# The commodity label lives in cell C4 of every yearly worksheet (2000-2100),
# all sharing the same underlying string, so update it everywhere it appears.
for ($i = 1; $i -le $wb.Worksheets.Count; $i++) {
    $ws = $wb.Worksheets.Item($i)
    if ($ws.Range("C4").Value2 -eq "Copper ores and concentrates") {
        $ws.Range("C4").Value = "Copper"
    }
}

# Small last-digit precision refinements to the D4 totals that accompanied
# the disaggregated source data, on the specific year sheets affected.
$wb.Worksheets.Item("2021").Range("D4").Value = 27629.08234046596
$wb.Worksheets.Item("2023").Range("D4").Value = 45474.52846901826
$wb.Worksheets.Item("2025").Range("D4").Value = 52615.60445701829
$wb.Worksheets.Item("2028").Range("D4").Value = 76414.31688861702
$wb.Worksheets.Item("2031").Range("D4").Value = 102815.2737995718
$wb.Worksheets.Item("2041").Range("D4").Value = 413896.0450908013
$wb.Worksheets.Item("2044").Range("D4").Value = 824812.1633242127
$wb.Worksheets.Item("2048").Range("D4").Value = 1595482.848064659
$wb.Worksheets.Item("2054").Range("D4").Value = 1715249.128188553
$wb.Worksheets.Item("2072").Range("D4").Value = 1598520.870762428
$wb.Worksheets.Item("2077").Range("D4").Value = 1443658.30372603
$wb.Worksheets.Item("2092").Range("D4").Value = 1706284.654525028
